$wb = $excel.ActiveWorkbook
$wsAssoziationen = $wb.Worksheets.Item("Assoziationen")

# The association "Kaiser/Mandarin/Elefant/Reiter/Turm/Geschütz/Bauer -> Spielfigur"
# was mislabeled as "Ist abgeleitet von" (is derived from); correct it to
# "Ist Teil von" (is part of) to match the updated class diagram / associations.
$wsAssoziationen.Range("A10").Value2 = "Ist Teil von"
$wsAssoziationen.Range("A11").Value2 = "Ist Teil von"
$wsAssoziationen.Range("A12").Value2 = "Ist Teil von"
$wsAssoziationen.Range("A13").Value2 = "Ist Teil von"
$wsAssoziationen.Range("A14").Value2 = "Ist Teil von"
$wsAssoziationen.Range("A15").Value2 = "Ist Teil von"
$wsAssoziationen.Range("A16").Value2 = "Ist Teil von"

# Move the active selection on the "Assoziationen" sheet, reflecting where the
# author was working when the associations table was revised.
$wsAssoziationen.Activate()
$wsAssoziationen.Range("B18").Select()
